# Applies the "Proposta de desnormalizacao da BD" edit to the Cronograma
# worksheet: rewrites several activity descriptions and appends two new
# schedule rows (9 and 10), then updates the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing activity rows (column B) ---
$ws.Range("B3").Value = "Construir  Interfaces de cadastro e consulta Cliente  (Jaime) falta Deletar"
$ws.Range("B4").Value = "Construir  Interfaces de cadastro e consulta Animal (Valter)"
$ws.Range("B5").Value = "Construir  Interfaces dos servicos Vacinacao"
$ws.Range("B6").Value = "Construir  Interfaces dos servicos Exames "
$ws.Range("B7").Value = "Construir  Interfaces dos servicos cirurgia"
$ws.Range("B8").Value = "Construir  Interfaces do historico"
$ws.Range("B9").Value = "Fazer relatorios"

# --- Append two new rows to the schedule ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Fazer a insercao, actualizacao e eliminacao  de varios dados "

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Tela de login (Senha) "

# --- Update selected/active cell ---
[void]$ws.Range("B4").Select()
